$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Row 65: Contest 53 (CSK vs PBKS) ----
# Raw score inputs (columns E,H,K,N,Q,T,W,Z,AC)
$ws.Range("E65").Value = 0
$ws.Range("H65").Value = 30
$ws.Range("K65").Value = 70
$ws.Range("N65").Value = 60
$ws.Range("Q65").Value = 50
$ws.Range("T65").Value = 0
$ws.Range("W65").Value = 100
$ws.Range("Z65").Value = 40
$ws.Range("AC65").Value = 80

# D65 and S65 tie for last place (8th/9th) so the formula is replaced with
# the manually-averaged value instead of the RANK/VLOOKUP formula.
$ws.Range("D65").Value = -22.5
$ws.Range("S65").Value = -22.5

# ---- Row 66: Contest 54 (KKR vs RR) ----
# Raw score inputs (columns E,H,K,N,Q,T,W,Z,AC)
$ws.Range("E66").Value = 30
$ws.Range("H66").Value = 100
$ws.Range("K66").Value = 70
$ws.Range("N66").Value = 80
$ws.Range("Q66").Value = 40
$ws.Range("T66").Value = 30
$ws.Range("W66").Value = 60
$ws.Range("Z66").Value = 50
$ws.Range("AC66").Value = 30

# D66, S66 and AB66 tie for last place (7th/8th/9th) so the formula is
# replaced with the manually-averaged value instead of the RANK/VLOOKUP formula.
$ws.Range("D66").Value = -20
$ws.Range("S66").Value = -20
$ws.Range("AB66").Value = -20
